$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (shifts old C..F to D..G) and inherits
# the adjacent (column B) formatting, matching the target style (s="2").
$ws.Range("C1").EntireColumn.Insert()

# Give the new PROJECT_ID column its own width (closest achievable value to
# the target 14.5234375 given the host's pixel-rounding on ColumnWidth).
$ws.Columns.Item(3).ColumnWidth = 13.65

# Populate the header and data for the new PROJECT_ID column. Shared strings
# get appended in first-use order: PROJECT_ID, Candy_kingdom, Nightosphere,
# Treehouse.
$ws.Range("C1").Value = "PROJECT_ID"
$ws.Range("C2").Value = "Candy_kingdom"
$ws.Range("C3").Value = "Nightosphere"
$ws.Range("C4").Value = "Treehouse"
$ws.Range("C5").Value = "Treehouse"
$ws.Range("C6").Value = "Treehouse"
$ws.Range("C7").Value = "Treehouse"
$ws.Range("C8").Value = "Treehouse"
$ws.Range("C9").Value = "Treehouse"
$ws.Range("C10").Value = "Treehouse"
$ws.Range("C11").Value = "Nightosphere"
$ws.Range("C12").Value = "Nightosphere"
$ws.Range("C13").Value = "Candy_kingdom"
$ws.Range("C14").Value = "Nightosphere"
$ws.Range("C15").Value = "Nightosphere"
$ws.Range("C16").Value = "Nightosphere"
$ws.Range("C17").Value = "Nightosphere"
$ws.Range("C18").Value = "Nightosphere"
$ws.Range("C19").Value = "Candy_kingdom"
$ws.Range("C20").Value = "Candy_kingdom"

# Match the author's final selection (whole PROJECT_ID column highlighted).
[void]$ws.Columns.Item(3).Select()
